# Automatische test-sync: 2025-07-27 19:48:50
# Append the newest test-mail row to the "Logs" sheet and bump the matching
# category tally on the "Dashboard" sheet. Also extends the Logs
# conditional-formatting ranges and the Dashboard bar-chart series ranges
# so they keep covering the newly added rows.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs!A19:J19 -----------------------------------------------------
$logsRow = 19

$logs.Cells.Item($logsRow, 1).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($logsRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($logsRow, 3).Value = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($logsRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($logsRow, 5).Value = "Beste,`nBedankt voor je e-mail. Wij zullen proberen om de demo in te plannen bij Van Dijk op vrijdag om 11:00 uur. `nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($logsRow, 6).Value = "2025-07-27 19:48:00"
$logs.Cells.Item($logsRow, 7).Value = "Ja"
$logs.Cells.Item($logsRow, 8).Value = "Nee"
$logs.Cells.Item($logsRow, 9).Value = "Ja"
$logs.Cells.Item($logsRow, 10).Value = "Nee"

# The multi-line "Antwoord" text auto-expands the row height; re-fit it back
# down to the sheet's normal single-line height (matches the other rows).
$logs.Rows.Item($logsRow).AutoFit()

# --- Dashboard!A7:B7 ----------------------------------------------------
$dashRow = 7

$dashboard.Cells.Item($dashRow, 1).Value = "Planning / Afspraak"
$dashboard.Cells.Item($dashRow, 2).Value = 1

# --- Logs: grow the conditional-formatting ranges to include row 19 -----
$cfCols = "D", "G", "H", "I", "J"
foreach ($col in $cfCols) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "18")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "19")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard chart: extend the category/value series ranges to row 7 --
$chart = $dashboard.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$7,Dashboard!`$B`$2:`$B`$7,1)"
